# Auto-generated edit script applying the Adamantoise_Profits diff
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.
$wb = $excel.ActiveWorkbook


$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 759.8
$ws.Range("I9").Value = 799.75
$ws.Range("J9").Value = 600
$ws.Range("K9").Value = 799.75
$ws.Range("L9").Value = 600
$ws.Range("M9").Value = -630.75
$ws.Range("N9").Value = -938
$ws.Range("H51").Value = 3234.5925
$ws.Range("I51").Value = 3407.913
$ws.Range("J51").Value = 2238
$ws.Range("K51").Value = 3407.913
$ws.Range("L51").Value = 2238
$ws.Range("M51").Value = -2923.913
$ws.Range("N51").Value = -3206
$ws.Range("H106").Value = 8335698
$ws.Range("J106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("N106").ClearContents()
$ws.Range("H137").Value = 3090279.8
$ws.Range("I137").Value = 2859.75
$ws.Range("J137").Value = 5560215.5
$ws.Range("K137").Value = 8579.25
$ws.Range("L137").Value = 16680646.5
$ws.Range("M137").Value = -6029.25
$ws.Range("N137").Value = -16685746.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 4199.3125
$ws.Range("I61").Value = 4400
$ws.Range("J61").Value = 3941.2856
$ws.Range("K61").Value = 4400
$ws.Range("L61").Value = 3941.2856
$ws.Range("M61").Value = -4188
$ws.Range("N61").Value = -4365.2856
$ws.Range("H104").Value = 106000
$ws.Range("J104").Value = 106000
$ws.Range("L104").Value = 106000
$ws.Range("N104").Value = -112988
$ws.Range("H132").Value = 2081.2058
$ws.Range("I132").Value = 2167.5334
$ws.Range("K132").Value = 6502.600199999999
$ws.Range("M132").Value = -3972.600199999999
$ws.Range("H136").Value = 4199.3125
$ws.Range("I136").Value = 4400
$ws.Range("J136").Value = 3941.2856
$ws.Range("K136").Value = 13200
$ws.Range("L136").Value = 11823.8568
$ws.Range("M136").Value = -10650
$ws.Range("N136").Value = -16923.8568

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2929.1738
$ws.Range("I86").Value = 2763.1765
$ws.Range("K86").Value = 2763.1765
$ws.Range("M86").Value = -1640.1765
$ws.Range("H89").Value = 2929.1738
$ws.Range("I89").Value = 2763.1765
$ws.Range("K89").Value = 13815.8825
$ws.Range("M89").Value = -8199.8825
$ws.Range("H92").Value = 55728.047
$ws.Range("J92").Value = 55728.047
$ws.Range("L92").Value = 55728.047
$ws.Range("N92").Value = -60720.047
$ws.Range("H94").Value = 619.913
$ws.Range("I94").Value = 693.8333
$ws.Range("J94").Value = 539.2727
$ws.Range("K94").Value = 693.8333
$ws.Range("L94").Value = 539.2727
$ws.Range("M94").Value = -242.8333
$ws.Range("N94").Value = -1441.2727
$ws.Range("H107").Value = 2472.1
$ws.Range("I107").Value = 3218.6155
$ws.Range("J107").Value = 1085.7142
$ws.Range("K107").Value = 3218.6155
$ws.Range("L107").Value = 1085.7142
$ws.Range("M107").Value = -1298.6155
$ws.Range("N107").Value = -4925.7142
$ws.Range("H132").Value = 76666.336
$ws.Range("J132").Value = 76666.336
$ws.Range("L132").Value = 76666.336
$ws.Range("N132").Value = -86786.336
$ws.Range("H134").Value = 1711406.4
$ws.Range("I134").Value = 1906715.8
$ws.Range("K134").Value = 5720147.4
$ws.Range("M134").Value = -5717612.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 69.666664
$ws.Range("I7").Value = 34.3
$ws.Range("J7").Value = 246.5
$ws.Range("K7").Value = 34.3
$ws.Range("L7").Value = 246.5
$ws.Range("M7").Value = 78.7
$ws.Range("N7").Value = -472.5
$ws.Range("H58").Value = 3112.4314
$ws.Range("I58").Value = 2670.8108
$ws.Range("K58").Value = 2670.8108
$ws.Range("M58").Value = -2467.8108
$ws.Range("H134").Value = 2694.16
$ws.Range("I134").Value = 2635.4546
$ws.Range("J134").Value = 3124.6667
$ws.Range("K134").Value = 7906.3638
$ws.Range("L134").Value = 9374.000100000001
$ws.Range("M134").Value = -5371.3638
$ws.Range("N134").Value = -14444.0001
$ws.Range("H136").Value = 3112.4314
$ws.Range("I136").Value = 2670.8108
$ws.Range("K136").Value = 8012.432400000001
$ws.Range("M136").Value = -5462.432400000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 6359.5835
$ws.Range("I3").Value = 6359.5835
$ws.Range("K3").Value = 19078.7505
$ws.Range("M3").Value = -18966.7505
$ws.Range("H5").Value = 822.75
$ws.Range("I5").Value = 698
$ws.Range("J5").Value = 897.6
$ws.Range("K5").Value = 2094
$ws.Range("L5").Value = 2692.8
$ws.Range("M5").Value = -1982
$ws.Range("N5").Value = -2916.8
$ws.Range("H87").Value = 1500
$ws.Range("I87").Value = 1500
$ws.Range("K87").Value = 4500
$ws.Range("M87").Value = -3252
$ws.Range("H90").Value = 1500
$ws.Range("I90").Value = 1500
$ws.Range("K90").Value = 13500
$ws.Range("M90").Value = -7260
$ws.Range("H107").Value = 722.0952
$ws.Range("J107").Value = 708.5
$ws.Range("L107").Value = 2125.5
$ws.Range("N107").Value = -5965.5
$ws.Range("H113").Value = 2165.875
$ws.Range("J113").Value = 2238
$ws.Range("L113").Value = 6714
$ws.Range("N113").Value = -11054
$ws.Range("H120").Value = 16413.37
$ws.Range("I120").Value = 6872.6665
$ws.Range("K120").Value = 20617.9995
$ws.Range("M120").Value = -15779.9995
$ws.Range("H122").Value = 941.8570999999999
$ws.Range("I122").Value = 1149.5
$ws.Range("K122").Value = 10345.5
$ws.Range("M122").Value = -7895.5
$ws.Range("H127").Value = 2445
$ws.Range("J127").Value = 2445
$ws.Range("L127").Value = 7335
$ws.Range("N127").Value = -17255
$ws.Range("H131").Value = 1622.3611
$ws.Range("I131").Value = 1170.4286
$ws.Range("J131").Value = 1731.4482
$ws.Range("K131").Value = 3511.2858
$ws.Range("L131").Value = 5194.3446
$ws.Range("M131").Value = 1528.7142
$ws.Range("N131").Value = -15274.3446
$ws.Range("H132").Value = 3336666.2
$ws.Range("I132").Value = 5000
$ws.Range("K132").Value = 45000
$ws.Range("M132").Value = -42470
$ws.Range("H134").Value = 2754.0833
$ws.Range("I134").Value = 2754.0833
$ws.Range("K134").Value = 8262.249899999999
$ws.Range("M134").Value = -3192.249899999999
$ws.Range("H135").Value = 822.75
$ws.Range("I135").Value = 698
$ws.Range("J135").Value = 897.6
$ws.Range("K135").Value = 6282
$ws.Range("L135").Value = 8078.400000000001
$ws.Range("M135").Value = -3747
$ws.Range("N135").Value = -13148.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H10").Value = 10001000
$ws.Range("H132").Value = 2415.0908
$ws.Range("I132").Value = 2652.2856
$ws.Range("K132").Value = 7956.8568
$ws.Range("M132").Value = -5426.8568
$ws.Range("H136").Value = 67506.92
$ws.Range("J136").Value = 67506.92
$ws.Range("L136").Value = 202520.76
$ws.Range("N136").Value = -207620.76

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2485.0833
$ws.Range("I7").Value = 2332.1
$ws.Range("K7").Value = 2332.1
$ws.Range("M7").Value = -2220.1
$ws.Range("H100").Value = 3570.2856
$ws.Range("I100").Value = 3998
$ws.Range("K100").Value = 3998
$ws.Range("M100").Value = -3457
$ws.Range("H126").Value = 2485.0833
$ws.Range("I126").Value = 2332.1
$ws.Range("K126").Value = 6996.299999999999
$ws.Range("M126").Value = -4526.299999999999
$ws.Range("H136").Value = 9461.956
$ws.Range("I136").Value = 9508.538
$ws.Range("J136").Value = 9401.4
$ws.Range("K136").Value = 28525.614
$ws.Range("L136").Value = 28204.2
$ws.Range("M136").Value = -25975.614
$ws.Range("N136").Value = -33304.2
$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").ClearContents()
$ws.Range("N141").Value = 0

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 12018.272
$ws.Range("I96").Value = 6690.6
$ws.Range("K96").Value = 6690.6
$ws.Range("M96").Value = -5317.6
$ws.Range("H123").Value = 61994
$ws.Range("J123").Value = 61994
$ws.Range("L123").Value = 61994
$ws.Range("N123").Value = -71794
$ws.Range("H126").Value = 3495.5
$ws.Range("I126").Value = 3495.5
$ws.Range("K126").Value = 10486.5
$ws.Range("M126").Value = -8016.5
$ws.Range("H132").Value = 3671.5
$ws.Range("I132").Value = 3671.5
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 11014.5
$ws.Range("L132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -8484.5
$ws.Range("H136").Value = 2715.389
$ws.Range("I136").Value = 2429.2307
$ws.Range("J136").Value = 3459.4
$ws.Range("K136").Value = 7287.6921
$ws.Range("L136").Value = 10378.2
$ws.Range("M136").Value = -4737.6921
$ws.Range("N136").Value = -15478.2
